$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("dnasr281@gmail.com, ")) {
        $rest = $val.Substring(20)
        $newVal = "$rest, dnasr281@gmail.com"
        $cell.Value2 = $newVal
    }
}
